$wb = $excel.ActiveWorkbook

# Update the Date value on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-03-16T01:14:24+00:00"

# Fix the EXOMIZER -> EXOMISER typo on the Concepts sheet
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("B8").Value = "EXOMISER"
$concepts.Range("C8").Value = "Exomiser Report"
